$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.929.66'
$ws.Range("E2").Value = '  +0.86%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.511.11'
$ws.Range("E3").Value = '  -0.34%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.44'
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.16'
$ws.Range("E6").Value = '  +4.26%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.594'
$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.138'
$ws.Range("E9").Value = '  +2.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.12'
$ws.Range("E10").Value = '  -2.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.433'
$ws.Range("E11").Value = '  -1.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.105.71'
$ws.Range("E12").Value = '  -0.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.21'
$ws.Range("E13").Value = '  +11.78%  '

$ws.Range("E14").Value = '  -0.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.816.65'
$ws.Range("E15").Value = '  +0.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000181'
$ws.Range("E16").Value = '  -1.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.506.42'
$ws.Range("E17").Value = '  +0.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.36'
$ws.Range("E18").Value = '  +0.26%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.72'
$ws.Range("E19").Value = '  +3.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '395.19'
$ws.Range("E20").Value = '  -0.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.02'
$ws.Range("E21").Value = '  +0.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.63'
$ws.Range("E22").Value = '  +0.21%  '

$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.543'
$ws.Range("E23").Value = '  +0.59%  '

$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.36%  '

$ws.Range("E25").Value = '  +0.46%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000124'
$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.38'
$ws.Range("E27").Value = '  +0.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.179'
$ws.Range("E28").Value = '  -1.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.996'
$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.27'
$ws.Range("E30").Value = '  -0.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.46'
$ws.Range("E31").Value = '  -0.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.07'
$ws.Range("E32").Value = '  -0.66%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.92'
$ws.Range("E33").Value = '  -1.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.41'
$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.66'
$ws.Range("E36").Value = '  +1.45%  '

$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.73'
$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.96'
$ws.Range("E38").Value = '  +2.27%  '

$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.876'
$ws.Range("E39").Value = '  -2.32%  '

$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.14'
$ws.Range("E40").Value = '  +2.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.73'
$ws.Range("E41").Value = '  +0.03%  '

$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.66'
$ws.Range("E42").Value = '  +0.38%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.70'
$ws.Range("E43").Value = '  +0.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.66'
$ws.Range("E44").Value = '  +1.30%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.829.58'
$ws.Range("E45").Value = '  +0.68%  '

$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0732'
$ws.Range("E46").Value = '  -1.91%  '

$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '42.40'
$ws.Range("E47").Value = '  -1.24%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0305'
$ws.Range("E48").Value = '  -0.99%  '

$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '344.92'
$ws.Range("E49").Value = '  +1.23%  '

$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.08'
$ws.Range("E50").Value = '  -1.59%  '

$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.65'
$ws.Range("E51").Value = '  +0.05%  '
